# Refresh the "Price" (D) and "Volume(1h)" (E) columns of the cryptos list
# with the latest scraped snapshot (GitHub Actions cron job).
#
# Every value on this sheet is stored as plain text (the upstream feed
# formats prices/percentages as strings, e.g. "2.544.44" uses "." as a
# thousands separator, which is not a valid Excel number anyway). Some of
# the new price strings, though, *do* look like plain decimals (e.g.
# "143.49"), and Excel's Range.Value setter auto-coerces those into real
# Number cells. To keep every refreshed cell a Text cell (matching the
# rest of the sheet / the source feed), such values are written with the
# cell pinned to the "Text" (`@`) number format, then the style is put back
# to "Normal" so no stray formatting sticks around.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = 'D2'; Text = '63.047.26' }
    @{ Cell = 'E2'; Text = '  -0.14%  ' }
    @{ Cell = 'D3'; Text = '2.560.74' }
    @{ Cell = 'E3'; Text = '  -0.28%  ' }
    @{ Cell = 'D4'; Text = '0.997' }
    @{ Cell = 'E4'; Text = '  -0.23%  ' }
    @{ Cell = 'D5'; Text = '586.25' }
    @{ Cell = 'E5'; Text = '  -0.14%  ' }
    @{ Cell = 'D6'; Text = '143.49' }
    @{ Cell = 'E6'; Text = '  -3.31%  ' }
    @{ Cell = 'D7'; Text = '0.998' }
    @{ Cell = 'E7'; Text = '  -0.14%  ' }
    @{ Cell = 'D8'; Text = '0.585' }
    @{ Cell = 'E8'; Text = '  -2.09%  ' }
    @{ Cell = 'D9'; Text = '0.106' }
    @{ Cell = 'E9'; Text = '  -2.84%  ' }
    @{ Cell = 'D10'; Text = '5.59' }
    @{ Cell = 'E10'; Text = '  -1.16%  ' }
    @{ Cell = 'E11'; Text = '  -0.20%  ' }
    @{ Cell = 'D12'; Text = '0.348' }
    @{ Cell = 'E12'; Text = '  -2.45%  ' }
    @{ Cell = 'D13'; Text = '27.05' }
    @{ Cell = 'E13'; Text = '  -1.90%  ' }
    @{ Cell = 'D14'; Text = '3.009.28' }
    @{ Cell = 'E14'; Text = '  -0.71%  ' }
    @{ Cell = 'D15'; Text = '62.882.15' }
    @{ Cell = 'E15'; Text = '  -0.20%  ' }
    @{ Cell = 'D16'; Text = '0.0000145' }
    @{ Cell = 'E16'; Text = '  -1.53%  ' }
    @{ Cell = 'D17'; Text = '2.544.95' }
    @{ Cell = 'E17'; Text = '  -0.67%  ' }
    @{ Cell = 'D18'; Text = '11.04' }
    @{ Cell = 'E18'; Text = '  -2.91%  ' }
    @{ Cell = 'D19'; Text = '340.23' }
    @{ Cell = 'E19'; Text = '  -1.12%  ' }
    @{ Cell = 'D20'; Text = '4.29' }
    @{ Cell = 'E20'; Text = '  -3.28%  ' }
    @{ Cell = 'D21'; Text = '6.59' }
    @{ Cell = 'E21'; Text = '  -3.99%  ' }
    @{ Cell = 'E22'; Text = '  +0.04%  ' }
    @{ Cell = 'D23'; Text = '67.20' }
    @{ Cell = 'E23'; Text = '  +0.92%  ' }
    @{ Cell = 'D24'; Text = '1.57' }
    @{ Cell = 'E24'; Text = '  +5.38%  ' }
    @{ Cell = 'D25'; Text = '1.60' }
    @{ Cell = 'E25'; Text = '  -1.66%  ' }
    @{ Cell = 'E26'; Text = '  -3.99%  ' }
    @{ Cell = 'E27'; Text = '  -0.04%  ' }
    @{ Cell = 'D28'; Text = '7.93' }
    @{ Cell = 'E28'; Text = '  -3.82%  ' }
    @{ Cell = 'D29'; Text = '8.17' }
    @{ Cell = 'E29'; Text = '  -3.27%  ' }
    @{ Cell = 'D30'; Text = '1.94' }
    @{ Cell = 'E30'; Text = '  -2.67%  ' }
    @{ Cell = 'D31'; Text = '468.23' }
    @{ Cell = 'E31'; Text = '  +1.29%  ' }
    @{ Cell = 'D32'; Text = '0.0₃0795' }
    @{ Cell = 'E32'; Text = '  -3.89%  ' }
    @{ Cell = 'D33'; Text = '1.66' }
    @{ Cell = 'E33'; Text = '  +2.86%  ' }
    @{ Cell = 'D34'; Text = '175.83' }
    @{ Cell = 'E34'; Text = '  -0.22%  ' }
    @{ Cell = 'D35'; Text = '0.997' }
    @{ Cell = 'E35'; Text = '  -0.27%  ' }
    @{ Cell = 'E36'; Text = '  -2.53%  ' }
    @{ Cell = 'D37'; Text = '18.76' }
    @{ Cell = 'E37'; Text = '  -2.34%  ' }
    @{ Cell = 'D38'; Text = '4.51' }
    @{ Cell = 'E38'; Text = '  -2.28%  ' }
    @{ Cell = 'E39'; Text = '  -0.12%  ' }
    @{ Cell = 'D40'; Text = '1.71' }
    @{ Cell = 'E40'; Text = '  -2.32%  ' }
    @{ Cell = 'D41'; Text = '40.05' }
    @{ Cell = 'E41'; Text = '  +1.50%  ' }
    @{ Cell = 'D42'; Text = '157.93' }
    @{ Cell = 'E42'; Text = '  +4.37%  ' }
    @{ Cell = 'D43'; Text = '3.69' }
    @{ Cell = 'E43'; Text = '  -3.83%  ' }
    @{ Cell = 'D44'; Text = '21.37' }
    @{ Cell = 'E44'; Text = '  +1.59%  ' }
    @{ Cell = 'D45'; Text = '0.626' }
    @{ Cell = 'E45'; Text = '  +1.99%  ' }
    @{ Cell = 'D46'; Text = '0.0533' }
    @{ Cell = 'E46'; Text = '  -3.09%  ' }
    @{ Cell = 'D47'; Text = '0.0959' }
    @{ Cell = 'E47'; Text = '  -1.83%  ' }
    @{ Cell = 'D48'; Text = '0.0235' }
    @{ Cell = 'E48'; Text = '  -2.25%  ' }
    @{ Cell = 'D49'; Text = '17.94' }
    @{ Cell = 'E49'; Text = '  -2.79%  ' }
    @{ Cell = 'E50'; Text = '  -0.33%  ' }
    @{ Cell = 'D51'; Text = '1.67' }
    @{ Cell = 'E51'; Text = '  -4.35%  ' }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    if ($u.Text -match '^[+-]?\d+(\.\d+)?$') {
        # Plain-decimal text (e.g. "143.49") - pin to Text format first so
        # it is not silently turned into a Number cell.
        $cell.NumberFormat = "@"
        $cell.Value = $u.Text
        $cell.Style = "Normal"
    } else {
        # Already unambiguous as text (multiple dots, %, spaces, subscripts, ...).
        $cell.Value = $u.Text
    }
}
